# Add a new "2022-Q4" sheet (holdings detail) right before the existing
# "2022-Q3" sheet, and add a corresponding summary row on the "总计" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Create the new "2022-Q4" detail sheet by duplicating "2022-Q3" (so
#    it inherits identical column layout / styles / borders), inserting
#    it immediately before "2022-Q3", then overwrite its values with the
#    2022-Q4 figures.
#
#    Note: columns D/E/F/G on these sheets hold numeric-looking values
#    stored as *text* (the source data keeps two decimal places etc.),
#    so a leading apostrophe is used to stop them being auto-converted
#    to numbers - same for the "513120" fund code in B3.
# ---------------------------------------------------------------------
$wsQ3 = $wb.Worksheets.Item("2022-Q3")
$wsQ3.Copy($wsQ3)

# The freshly created copy is placed immediately before "2022-Q3" and
# becomes the active sheet.
$wsQ4 = $wb.ActiveSheet
$wsQ4.Name = "2022-Q4"

# Row 2: 010010 fund keeps its code/name, only the metrics change.
$wsQ4.Range("D2").Value = "'7.94"
$wsQ4.Range("E2").Value = "'97.04"
$wsQ4.Range("F2").Value = "'3.51"
$wsQ4.Range("G2").Value = "'0.2787"
$wsQ4.Range("H2").Value = 8

# Row 3: fund changes entirely (005646 -> 513120).
$wsQ4.Range("B3").Value = "'513120"
$wsQ4.Range("C3").Value = "广发中证香港创新药（QDII-ETF）"
$wsQ4.Range("D3").Value = "'1.85"
$wsQ4.Range("E3").Value = "'98.69"
$wsQ4.Range("F3").Value = "'2.71"
$wsQ4.Range("G3").Value = "'0.0501"
$wsQ4.Range("H3").Value = 10

# ---------------------------------------------------------------------
# 2) Update the "总计" (summary) sheet: insert a new row 2 for 2022-Q4,
#    pushing the existing quarters down by one row.
# ---------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("总计")
$wsTotal.Rows.Item(2).Insert()

# Restore formatting on the newly-inserted row by cloning it from the
# row right below (which still carries the original style).
$wsTotal.Range("A3:D3").Copy($wsTotal.Range("A2:D2"))

$wsTotal.Range("B2").Value = "2022-Q4"
$wsTotal.Range("C2").Value = 2
$wsTotal.Range("D2").Value = 0.33

# The "A" column is just a plain 0-based row counter independent of the
# shift above; re-number it for every data row (2..9) to keep it tidy.
for ($i = 0; $i -le 7; $i++) {
    $row = 2 + $i
    $wsTotal.Range("A$row").Value = $i
}
